$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 204, pushing existing rows 204..264 down to 205..265
$ws.Rows.Item(204).Insert()

# Populate the newly inserted row 204 with the new record
$ws.Range("A204").Value = 5
$ws.Range("B204").Value = "Macroferia Regional de Talca"
$ws.Range("C204").Value = "Maule"
$ws.Range("D204").Value = 44736
$ws.Range("D204").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E204").Value = 7
$ws.Range("F204").Value = 100112009
$ws.Range("G204").Value = "Acelga"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 300
$ws.Range("K204").Value = 3000
$ws.Range("L204").Value = 3000
$ws.Range("M204").Value = 3000
$ws.Range("N204").Value = "$/docena de atados (4 kilos)"
$ws.Range("O204").Value = "Región del Maule"
$ws.Range("P204").Value = 750
$ws.Range("Q204").Value = 4
$ws.Range("R204").Value = "Hortaliza"
